# Apply the commit: for each data row (2 through 40) on the active sheet,
# decrement the birth_year value in column Q by 1 and increment the
# age_y value in column S by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 40; $row++) {
    $qCell = $ws.Cells.Item($row, 17)   # column Q = birth_year
    $sCell = $ws.Cells.Item($row, 19)   # column S = age_y

    $qVal = $qCell.Value()
    $sVal = $sCell.Value()

    $qCell.Value = $qVal - 1
    $sCell.Value = $sVal + 1
}
